$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed coefficient / coef_pos values across rows 2-29
$ws.Range("H2").Value = 0.09433620093527373
$ws.Range("B3").Value = 0.1029574324638085
$ws.Range("H3").Value = 0.1972936333990822
$ws.Range("B4").Value = 0.08409957213983206
$ws.Range("H4").Value = 0.1784357730751058
$ws.Range("B5").Value = 0.06457353709297113
$ws.Range("H5").Value = 0.1589097380282449
$ws.Range("B6").Value = 0.04573420860271862
$ws.Range("H6").Value = 0.1400704095379924
$ws.Range("B7").Value = 0.02323290200375484
$ws.Range("C7").Value = 0.002508239712143095
$ws.Range("D7").Value = 3.018807889841133
$ws.Range("E7").Value = 0.0137343863590084
$ws.Range("F7").Value = 0.01830203432820247
$ws.Range("G7").Value = 0.02816376967930696
$ws.Range("H7").Value = 0.1175691029390286
$ws.Range("B8").Value = 0.02301880827754807
$ws.Range("C8").Value = 0.002391493838836631
$ws.Range("D8").Value = 3.962271444860268
$ws.Range("E8").Value = 0.00751086606513807
$ws.Range("F8").Value = 0.01830931334716203
$ws.Range("G8").Value = 0.02772830320793403
$ws.Range("H8").Value = 0.1173550092128218
$ws.Range("B9").Value = 0.02250608706600316
$ws.Range("H9").Value = 0.1168422880012769
$ws.Range("B10").Value = 0.02118775533242945
$ws.Range("H10").Value = 0.1155239562677032
$ws.Range("B11").Value = 0.03290212680018208
$ws.Range("H11").Value = 0.1272383277354558
$ws.Range("B12").Value = 0.05522884478585557
$ws.Range("H12").Value = 0.1495650457211293
$ws.Range("B13").Value = 0.06447114349309434
$ws.Range("H13").Value = 0.1588073444283681
$ws.Range("B14").Value = 0.07477726774693204
$ws.Range("H14").Value = 0.1691134686822058
$ws.Range("B15").Value = 0.08379256057805502
$ws.Range("H15").Value = 0.1781287615133287
$ws.Range("B16").Value = 0.08547196115689153
$ws.Range("H16").Value = 0.1798081620921653
$ws.Range("B17").Value = 0.08647074735830616
$ws.Range("H17").Value = 0.1808069482935799
$ws.Range("B18").Value = -0.09433620093527373
$ws.Range("B19").Value = 0.08991652577803222
$ws.Range("H19").Value = 0.1842527267133059
$ws.Range("B20").Value = 0.09289966074206668
$ws.Range("H20").Value = 0.1872358616773404
$ws.Range("B21").Value = 0.09690252412819515
$ws.Range("H21").Value = 0.1912387250634689
$ws.Range("B22").Value = 0.1024495846870818
$ws.Range("H22").Value = 0.1967857856223555
$ws.Range("B23").Value = 0.1066218907069758
$ws.Range("H23").Value = 0.2009580916422495
$ws.Range("B24").Value = 0.1085586091990521
$ws.Range("C24").Value = 0.009358680330341335
$ws.Range("D24").Value = 1902149455236.314
$ws.Range("E24").Value = 0.05414537341580854
$ws.Range("F24").Value = 0.09013975715596718
$ws.Range("G24").Value = 0.1269774612421372
$ws.Range("H24").Value = 0.2028948101343259
$ws.Range("B25").Value = 0.1104705985162107
$ws.Range("C25").Value = 0.008747141398140494
$ws.Range("D25").Value = 1181190483800.336
$ws.Range("E25").Value = 0.05263272765295209
$ws.Range("F25").Value = 0.093262762906375
$ws.Range("G25").Value = 0.1276784341260466
$ws.Range("H25").Value = 0.2048067994514844
$ws.Range("B26").Value = 0.1132252464187097
$ws.Range("C26").Value = 0.008672562860085362
$ws.Range("D26").Value = 1933663733577.018
$ws.Range("E26").Value = 0.05438712232633147
$ws.Range("F26").Value = 0.09617233815460344
$ws.Range("G26").Value = 0.1302781546828162
$ws.Range("H26").Value = 0.2075614473539834
$ws.Range("B27").Value = 0.1154173152168604
$ws.Range("C27").Value = 0.008791485020300619
$ws.Range("D27").Value = 21.55636275087432
$ws.Range("E27").Value = 0.06375662027712767
$ws.Range("F27").Value = 0.09811727210318105
$ws.Range("G27").Value = 0.1327173583305405
$ws.Range("H27").Value = 0.2097535161521341
$ws.Range("B28").Value = 0.1168375201834414
$ws.Range("C28").Value = 0.008268580323377113
$ws.Range("D28").Value = 1628309554361.351
$ws.Range("E28").Value = 0.06936834759677808
$ws.Range("F28").Value = 0.100572693624265
$ws.Range("G28").Value = 0.1331023467426175
$ws.Range("H28").Value = 0.2111737211187151
$ws.Range("B29").Value = 0.02085980600167717
$ws.Range("C29").Value = 0.002164877931745044
$ws.Range("D29").Value = 2.025833849636871
$ws.Range("E29").Value = 0.01876538776019561
$ws.Range("F29").Value = 0.0165941283720794
$ws.Range("G29").Value = 0.02512548363127512
$ws.Range("H29").Value = 0.1151960069369509

# Rows that lost their std-err/t/p/CI stats (now blank like other single-column rows)
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("G10").ClearContents()
